# Budget Simulation - Admin: add "Budget Simulation" / "Name Your Price" section
# headers, fix the escalator-row ordering, and add a duplicated
# "Name Your Price" bounds table (rows 17-19) below the existing one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Row 4 becomes the new "Budget Simulation " merged title row.
#    Seed fillId (theme fill used throughout this block) from C5's
#    format, then merge + apply bold/centered 14pt Calibri to the
#    whole row, but only bold the left (merged) cell so the look
#    matches the rest of the "VARIABLES" block.
# ---------------------------------------------------------------------
$ws.Range("C5").Copy()
$ws.Range("C4:E4").PasteSpecial($xlPasteFormats)
$ws.Range("C4:E4").Merge()
$ws.Range("C4:E4").Font.Name = "Calibri"
$ws.Range("C4:E4").Font.Size = 14
$ws.Range("C4:E4").HorizontalAlignment = $xlCenter
$ws.Range("C4:E4").VerticalAlignment = $xlCenter
$ws.Cells.Item(4, 3).Font.Bold = $true
$ws.Cells.Item(4, 3).Value = "Budget Simulation "
$ws.Rows.Item(4).RowHeight = 18

# ---------------------------------------------------------------------
# 2) Shift the escalator rows down by one (old row N -> new row N+1),
#    fixing the Overhead/Other order at the same time, by writing the
#    target text/value directly - the cell styles (s23/s23/s15) are
#    already correct for rows 5-9, row 10 needs the format copied in.
# ---------------------------------------------------------------------
$ws.Range("C9:E9").Copy()
$ws.Range("C10:E10").PasteSpecial($xlPasteFormats)

$ws.Range("C5").Value = "Contract Escalator "
$ws.Range("D5").Value = "(%)"
$ws.Range("E5").Value = 0.03

$ws.Range("C6").Value = "Diesel Price Escalator "
$ws.Range("D6").Value = "(%)"
$ws.Range("E6").Value = 0.05

$ws.Range("C7").Value = "Overhead Cost Allocation"
$ws.Range("D7").Value = "(%)"
$ws.Range("E7").Value = 0.15

$ws.Range("C8").Value = "Other Costs Escalator"
$ws.Range("D8").Value = "(%)"
$ws.Range("E8").Value = 0.02

$ws.Range("C9").Value = "M&R Escalator First Half-Life"
$ws.Range("D9").Value = "(%)"
$ws.Range("E9").Value = 0.06

$ws.Range("C10").Value = "M&R Escalator Second Half-Life "
$ws.Range("D10").Value = "(%)"
$ws.Range("E10").Value = 0.08

# ---------------------------------------------------------------------
# 3) Old row 11 (Highland Contract) is dropped entirely.
# ---------------------------------------------------------------------
$ws.Range("C11:E11").Clear()
$ws.Rows.Item(11).RowHeight = 18

# ---------------------------------------------------------------------
# 4) Rows 12-14 keep their existing styling, shift up by one row and
#    pick up the Deployment/Contract-term bound values.
# ---------------------------------------------------------------------
$ws.Range("C12").Value = "Deployment Year Lower Bound"
$ws.Range("D12").Value = "(≥)"
$ws.Range("E12").Value = 2022

$ws.Range("C13").Value = "Contract Term Lower Bound"
$ws.Range("D13").Value = "(≥)"
$ws.Range("E13").Value = 5

$ws.Range("C14").Value = "Contract Term Upper Bound"
$ws.Range("D14").Value = "(≤)"
$ws.Range("E14").Value = 20

# ---------------------------------------------------------------------
# 5) Row 16 becomes the "Name Your Price " merged title row - same
#    look as row 4, but (per the source workbook) the bold/centered
#    font is applied across the whole merged range uniformly.
# ---------------------------------------------------------------------
$ws.Range("C5").Copy()
$ws.Range("C16:E16").PasteSpecial($xlPasteFormats)
$ws.Range("C16:E16").Merge()
$ws.Range("C16:E16").Font.Name = "Calibri"
$ws.Range("C16:E16").Font.Size = 14
$ws.Range("C16:E16").Font.Bold = $true
$ws.Range("C16:E16").HorizontalAlignment = $xlCenter
$ws.Range("C16:E16").VerticalAlignment = $xlCenter
$ws.Cells.Item(16, 3).Value = "Name Your Price "

# ---------------------------------------------------------------------
# 6) Rows 17-19: a duplicate "Name Your Price" bounds table, reusing
#    the Contract Term Lower/Upper Bound rows and adding the brand new
#    "Years to Deploy Upper Bound" row.
# ---------------------------------------------------------------------
$ws.Range("C13:E14").Copy()
$ws.Range("C17:E18").PasteSpecial($xlPasteFormats)
$ws.Range("C14:E14").Copy()
$ws.Range("C19:E19").PasteSpecial($xlPasteFormats)

$ws.Range("C17").Value = "Contract Term Lower Bound"
$ws.Range("D17").Value = "(≥)"
$ws.Range("E17").Value = 5

$ws.Range("C18").Value = "Contract Term Upper Bound"
$ws.Range("D18").Value = "(≤)"
$ws.Range("E18").Value = 20

$ws.Range("C19").Value = "Years to Deploy Upper Bound"
$ws.Range("D19").Value = "(≤)"
$ws.Range("E19").Value = 20

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 7) Selection, matching the saved cursor position in the workbook.
# ---------------------------------------------------------------------
[void]$ws.Range("I7").Select()
